# Rename the two sheets to unify the DataNode / DataTable / Entity
# terminology (formerly "Property1" / "Property2").
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$ws1.Name = "DataNode_1"
$ws2.Name = "DataNode_2"

# Row-1 / Row-8 header heights were re-flowed (autofit) when the sheet was
# last resaved; pin them to the values the workbook now carries.
$ws1.Rows.Item(1).RowHeight = 27
$ws1.Rows.Item(8).RowHeight = 40.5

$ws2.Rows.Item(1).RowHeight = 27
$ws2.Rows.Item(8).RowHeight = 67.5

# The second sheet ("DataNode_2") is now the active/visible tab.
$ws2.Activate()
